$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-05 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-06 Sunday", 2)

$d.Content.Find.Execute("378÷6=63, 0", $true, $false, $false, $false, $false, $true, 1, $false, "283÷7=40, 3", 2)
$d.Content.Find.Execute("518÷2=259, 0", $true, $false, $false, $false, $false, $true, 1, $false, "550÷4=137, 2", 2)
$d.Content.Find.Execute("719÷5=143, 4", $true, $false, $false, $false, $false, $true, 1, $false, "377÷4=94, 1", 2)
$d.Content.Find.Execute("801÷6=133, 3", $true, $false, $false, $false, $false, $true, 1, $false, "167÷2=83, 1", 2)
$d.Content.Find.Execute("103÷8=12, 7", $true, $false, $false, $false, $false, $true, 1, $false, "312÷7=44, 4", 2)

$d.Content.Find.Execute("516÷8=64, 4", $true, $false, $false, $false, $false, $true, 1, $false, "975÷3=325, 0", 2)
$d.Content.Find.Execute("129÷9=14, 3", $true, $false, $false, $false, $false, $true, 1, $false, "894÷9=99, 3", 2)
$d.Content.Find.Execute("415÷9=46, 1", $true, $false, $false, $false, $false, $true, 1, $false, "876÷9=97, 3", 2)
$d.Content.Find.Execute("434÷4=108, 2", $true, $false, $false, $false, $false, $true, 1, $false, "581÷7=83, 0", 2)
$d.Content.Find.Execute("757÷2=378, 1", $true, $false, $false, $false, $false, $true, 1, $false, "661÷8=82, 5", 2)

$d.Content.Find.Execute("582÷5=116, 2", $true, $false, $false, $false, $false, $true, 1, $false, "307÷2=153, 1", 2)
$d.Content.Find.Execute("554÷5=110, 4", $true, $false, $false, $false, $false, $true, 1, $false, "988÷3=329, 1", 2)
$d.Content.Find.Execute("654÷4=163, 2", $true, $false, $false, $false, $false, $true, 1, $false, "575÷6=95, 5", 2)
$d.Content.Find.Execute("151÷4=37, 3", $true, $false, $false, $false, $false, $true, 1, $false, "473÷4=118, 1", 2)
$d.Content.Find.Execute("582÷9=64, 6", $true, $false, $false, $false, $false, $true, 1, $false, "173÷9=19, 2", 2)

$d.Content.Find.Execute("303÷3=101, 0", $true, $false, $false, $false, $false, $true, 1, $false, "792÷8=99, 0", 2)
$d.Content.Find.Execute("593÷2=296, 1", $true, $false, $false, $false, $false, $true, 1, $false, "140÷7=20, 0", 2)
$d.Content.Find.Execute("954÷7=136, 2", $true, $false, $false, $false, $false, $true, 1, $false, "969÷4=242, 1", 2)
$d.Content.Find.Execute("412÷9=45, 7", $true, $false, $false, $false, $false, $true, 1, $false, "784÷9=87, 1", 2)
$d.Content.Find.Execute("316÷5=63, 1", $true, $false, $false, $false, $false, $true, 1, $false, "921÷8=115, 1", 2)

$d.Content.Find.Execute("926÷7=132, 2", $true, $false, $false, $false, $false, $true, 1, $false, "230÷2=115, 0", 2)
$d.Content.Find.Execute("368÷4=92, 0", $true, $false, $false, $false, $false, $true, 1, $false, "284÷9=31, 5", 2)
$d.Content.Find.Execute("744÷2=372, 0", $true, $false, $false, $false, $false, $true, 1, $false, "308÷6=51, 2", 2)
$d.Content.Find.Execute("389÷3=129, 2", $true, $false, $false, $false, $false, $true, 1, $false, "921÷2=460, 1", 2)
$d.Content.Find.Execute("643÷8=80, 3", $true, $false, $false, $false, $false, $true, 1, $false, "178÷7=25, 3", 2)
